$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.863.04"
$ws.Range("E2").Value = "  -1.84%  "
$ws.Range("D3").Value = "1.809.87"
$ws.Range("E3").Value = "  -0.86%  "
$ws.Range("E4").Value = "  +0.12%  "
$ws.Range("D5").Value = "'309.63"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.38%  "
$ws.Range("E6").Value = "  +0.11%  "
$ws.Range("D7").Value = "'0.4654"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.62%  "
$ws.Range("D8").Value = "'0.3696"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.49%  "
$ws.Range("D9").Value = "'0.07357"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -1.12%  "
$ws.Range("D10").Value = "'0.8726"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.45%  "
$ws.Range("D11").Value = "'20.46"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.62%  "
$ws.Range("D12").Value = "1.850.41"
$ws.Range("E12").Value = "  +1.34%  "
$ws.Range("D13").Value = "'5.354"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.46%  "
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'6.509"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.75%  "
$ws.Range("B15").Value = "TRON"
$ws.Range("C15").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D15").Value = "'0.07051"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.63%  "
$ws.Range("D16").Value = "'91.29"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.94%  "
$ws.Range("D17").Value = "'1.001"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.06%  "
$ws.Range("D18").Value = "'0.000008699"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.15%  "
$ws.Range("D19").Value = "'1.001"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.10%  "
$ws.Range("D20").Value = "'14.70"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.28%  "
$ws.Range("D21").Value = "26.878.80"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "'5.319"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.15%  "
$ws.Range("D23").Value = "'10.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.13%  "
$ws.Range("D24").Value = "2.060.17"
$ws.Range("E24").Value = "  +0.40%  "
$ws.Range("D25").Value = "'1.901"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.12%  "
$ws.Range("D26").Value = "'151.63"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.25%  "
$ws.Range("D27").Value = "'18.37"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "'2.142"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -5.23%  "
$ws.Range("D29").Value = "'5.305"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.68%  "
$ws.Range("D30").Value = "'115.86"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.10%  "
$ws.Range("D31").Value = "'0.08880"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.38%  "
$ws.Range("D32").Value = "'0.7528"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.85%  "
$ws.Range("D33").Value = "'1.153"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.54%  "
$ws.Range("D34").Value = "'2.917"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.61%  "
$ws.Range("D35").Value = "'4.456"
$ws.Range("D35").Style = "Normal"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").Value = "'0.01963"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.63%  "
$ws.Range("D39").Value = "'0.05252"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.02%  "
$ws.Range("D40").Value = "'2.429"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.69%  "
$ws.Range("D41").Value = "'2.926"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.08%  "
$ws.Range("D42").Value = "'0.5308"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.86%  "
$ws.Range("D43").Value = "'7.162"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.83%  "
$ws.Range("D44").Value = "'0.1663"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.39%  "
$ws.Range("D45").Value = "'8.446"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.44%  "
$ws.Range("D46").Value = "'0.4936"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -3.14%  "
$ws.Range("D47").Value = "'10.36"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -2.49%  "
$ws.Range("E48").Value = "  +0.13%  "
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("D50").Value = "'103.01"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.61%  "
$ws.Range("D51").Value = "'0.06279"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.63%  "
